$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mapping_com")

# Replace the "CR/LFM" prefix with "CR/LFINF" in the commercial (Offices, Trade,
# Hotels) mapping scheme blocks, mirroring the residential mapping scheme.
foreach ($addr in @("B2", "C2", "D2")) {
    $cell = $ws.Range($addr)
    $text = $cell.Value()
    $cell.Value = ($text -replace "CR/LFM", "CR/LFINF")
}
